$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 5599.6665
$ws.Range("I21").Value = 6620
$ws.Range("J21").Value = 498
$ws.Range("K21").Value = 6620
$ws.Range("L21").Value = 498
$ws.Range("M21").Value = -6152
$ws.Range("N21").Value = -1434

$ws.Range("H23").Value = 5599.6665
$ws.Range("I23").Value = 6620
$ws.Range("J23").Value = 498
$ws.Range("K23").Value = 6620
$ws.Range("L23").Value = 498
$ws.Range("M23").Value = -6386
$ws.Range("N23").Value = -966

$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3685
$ws.Range("N76").Value = -4630

$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2908
$ws.Range("N79").Value = -6184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H45").Value = 2194.353
$ws.Range("I45").Value = 1827.0834
$ws.Range("J45").Value = 3075.8
$ws.Range("K45").Value = 1827.0834
$ws.Range("L45").Value = 3075.8
$ws.Range("M45").Value = -1450.0834
$ws.Range("N45").Value = -3829.8

$ws.Range("H102").Value = 4388.9165
$ws.Range("I102").Value = 3266.7
$ws.Range("K102").Value = 3266.7
$ws.Range("M102").Value = -1644.7

$ws.Range("H135").Value = 32500
$ws.Range("J135").Value = 32500
$ws.Range("L135").Value = 32500
$ws.Range("N135").Value = -42640

$ws.Range("H139").Value = 20499.5
$ws.Range("J139").Value = 20499.5
$ws.Range("L139").Value = 20499.5
$ws.Range("N139").Value = -30779.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 25463.285
$ws.Range("J95").Value = 25463.285
$ws.Range("L95").Value = 25463.285
$ws.Range("N95").Value = -30955.285

$ws.Range("H135").Value = 97509
$ws.Range("J135").Value = 97509
$ws.Range("L135").Value = 97509
$ws.Range("N135").Value = -107649

$ws.Range("H137").Value = 75000
$ws.Range("I137").Value = 20000
$ws.Range("K137").Value = 20000
$ws.Range("M137").Value = -14900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1817.4286
$ws.Range("J23").Value = 1833.3334
$ws.Range("L23").Value = 1833.3334
$ws.Range("N23").Value = -2313.3334

$ws.Range("H27").Value = 1817.4286
$ws.Range("J27").Value = 1833.3334
$ws.Range("L27").Value = 1833.3334
$ws.Range("N27").Value = -2217.3334

$ws.Range("H39").Value = 908.6667
$ws.Range("I39").Value = 908.6667
$ws.Range("K39").Value = 908.6667
$ws.Range("M39").Value = -517.6667

$ws.Range("H49").Value = 908.6667
$ws.Range("I49").Value = 908.6667
$ws.Range("K49").Value = 908.6667
$ws.Range("M49").Value = -726.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 7486.6
$ws.Range("I62").Value = 100
$ws.Range("J62").Value = 9333.25
$ws.Range("K62").Value = 300
$ws.Range("L62").Value = 27999.75
$ws.Range("M62").Value = 386
$ws.Range("N62").Value = -29371.75

$ws.Range("H65").Value = 7486.6
$ws.Range("I65").Value = 100
$ws.Range("J65").Value = 9333.25
$ws.Range("K65").Value = 900
$ws.Range("L65").Value = 83999.25
$ws.Range("M65").Value = 2532
$ws.Range("N65").Value = -90863.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

$ws.Range("H132").Value = 4749.174
$ws.Range("J132").Value = 6638.5557
$ws.Range("L132").Value = 19915.6671
$ws.Range("N132").Value = -24975.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6794.5
$ws.Range("I68").Value = 2981.6667
$ws.Range("K68").Value = 2981.6667
$ws.Range("M68").Value = -2232.6667

$ws.Range("H71").Value = 6794.5
$ws.Range("I71").Value = 2981.6667
$ws.Range("K71").Value = 14908.3335
$ws.Range("M71").Value = -11164.3335

$ws.Range("H82").Value = 5775.25
$ws.Range("I82").Value = 2002
$ws.Range("J82").Value = 6314.2856
$ws.Range("K82").Value = 2002
$ws.Range("L82").Value = 6314.2856
$ws.Range("M82").Value = -1641
$ws.Range("N82").Value = -7036.2856

$ws.Range("H85").Value = 5775.25
$ws.Range("I85").Value = 2002
$ws.Range("J85").Value = 6314.2856
$ws.Range("K85").Value = 2002
$ws.Range("L85").Value = 6314.2856
$ws.Range("M85").Value = -754
$ws.Range("N85").Value = -8810.285599999999

$ws.Range("H93").Value = 2001.5
$ws.Range("I93").Value = 1003
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 1003
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = 245
$ws.Range("N93").Value = -5496

$ws.Range("H127").Value = 68943
$ws.Range("J127").Value = 68943
$ws.Range("L127").Value = 68943
$ws.Range("N127").Value = -78863

$ws.Range("H135").Value = 100214.5
$ws.Range("J135").Value = 100214.5
$ws.Range("L135").Value = 100214.5
$ws.Range("N135").Value = -110354.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11000
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 11000
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -66240

$ws.Range("H118").Value = 52896
$ws.Range("J118").Value = 52896
$ws.Range("L118").Value = 52896
$ws.Range("N118").Value = -56210

$ws.Range("H136").Value = 3145.6956
$ws.Range("I136").Value = 2290.0667
$ws.Range("K136").Value = 6870.2001
$ws.Range("M136").Value = -4320.2001
